# Move regression estimates to country specific folder
# Rename worksheet "EL" to "Student_share"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EL")
$ws.Name = "Student_share"
